# Importer: prevent product to be linked to variation child
# Insert a new test row before current row 46 (pushing the existing
# "variation/parent sku" negative-test rows down by one) and populate it
# with a new negative-test case: a variation whose parent SKU points to
# another variation child (row previously at 46 becomes 47, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Shift rows 46:50 down to 47:51 by inserting a blank row at 46.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new test case.
$ws.Cells.Item(46, 2).Value = 42
$ws.Cells.Item(46, 3).Value = 28
$ws.Cells.Item(46, 4).Value = "This tries to link variation to child"
$ws.Cells.Item(46, 6).Value = "Color/Black"
$ws.Cells.Item(46, 7).Value = "Size/XS"
$ws.Cells.Item(46, 9).Value = 12
$ws.Cells.Item(46, 10).Value = 100
$ws.Cells.Item(46, 11).Value = "Test Category"
$ws.Cells.Item(46, 12).Value = "Test Category"
$ws.Cells.Item(46, 13).Value = "shirt1.jpeg"
$ws.Cells.Item(46, 14).Value = "shirt2.jpeg,shirt3.jpeg"

# Match the saved selection state from the edit (cell B52 selected).
$ws.Range("B52").Select()
